# Generate Report for Handback
# Refresh the handoff/handback timestamps for the 14d7b9cd-bc5d-4809-a087-4877fc173bc5
# report row (row 2) on both locale sheets, leaving the f3289838 row (row 3) untouched.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-14 06:35:42"
$zhcn.Range("H2").Value = "2016-03-14 06:35:56"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-14 06:35:45"
$dede.Range("H2").Value = "2016-03-14 06:36:02"
